$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column H ("Save"), using same style as other headers (s="1")
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Save column data values, row 2 through row 24
$saveValues = @(1,0,1,0,1,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,1,1)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
